# Validacion examen Roque Estructura de Datos Dual
# Adds a new record row (01NOV / EDA / 02 / 04D / RRC / MARR / V_00001)
# to the bottom of the "Hoja2" validation table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# --- New row 22 --------------------------------------------------------
# Column B carries the delivery day, formatted like the other "d-mmm"
# style day cells already in the sheet; leading "'" forces text entry
# (quote-prefix) exactly like the existing "DíaEntregado" values.
$ws.Range("B22").NumberFormat = "d-mmm"
$ws.Range("B22").Value = "'01NOV"

$ws.Range("C22").Value = "EDA"

$ws.Range("D22").Value = "'02"

# Write F before E so the new shared-string table gets the same
# insertion order as the source workbook (…, 01NOV, RRC, 04D).
$ws.Range("F22").Value = "'RRC"
$ws.Range("E22").Value = "'04D"

$ws.Range("G22").Value = "MARR"
$ws.Range("H22").Value = "V_00001"

$ws.Range("I22").Formula = "=CONCATENATE(B22,C22,D22,E22,F22,G22,H22)"

# Match the saved view state (active cell moved below the new row).
[void]$ws.Range("H23").Select()
